$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate the new flow-cytometry "FlowContour" lookup table (G2:N10) ---
# Cell writes are ordered so the shared-string table is built up in the same
# sequence the original workbook used (first occurrence of each unique label).

# Header row (bold) + first-seen labels
$ws.Range("G2").Value = "Tube"
$ws.Range("H2").Value = "FSC"
$ws.Range("I2").Value = "SSC"
$ws.Range("J2").Value = "FL1"
$ws.Range("K2").Value = "FL2"
$ws.Range("L2").Value = "FL3"
$ws.Range("M2").Value = "FL4"
$ws.Range("N2").Value = "FL5"

$ws.Range("L3").Value = "CD45-ECD"

$ws.Range("J10").Value = "NS1"
$ws.Range("K10").Value = "NS2"
$ws.Range("L10").Value = "NS3"
$ws.Range("M10").Value = "NS4"
$ws.Range("N10").Value = "NS5"

$ws.Range("J3").Value = "IgG1-FITC"
$ws.Range("J4").Value = "Kappa-FIT"
$ws.Range("J5").Value = "CD7-FITC"
$ws.Range("J6").Value = "CD15-FITC"
$ws.Range("J7").Value = "CD14-FITC"
$ws.Range("J8").Value = "HLA-DR-FITC"
$ws.Range("J9").Value = "CD5-FITC"

$ws.Range("K4").Value = "Lambda-PE"
$ws.Range("K3").Value = "IgG1-PE"
$ws.Range("K5").Value = "CD4-PE"
$ws.Range("K6").Value = "CD13-PE"
$ws.Range("K7").Value = "CD11c-PE"
$ws.Range("K8").Value = "CD117-PE"
$ws.Range("K9").Value = "CD19-PE"

$ws.Range("M3").Value = "IgG1-PC5"
$ws.Range("M4").Value = "CD19-PC5"
$ws.Range("M5").Value = "CD8-PC5"
$ws.Range("M6").Value = "CD16-PC5"
$ws.Range("M7").Value = "CD64-PC5"
$ws.Range("M8").Value = "CD34-PC5"
$ws.Range("M9").Value = "CD3-PC5"

$ws.Range("N3").Value = "IgG1-PC7"
$ws.Range("N4").Value = "CD20-PC7"
$ws.Range("N5").Value = "CD2-PC7"
$ws.Range("N6").Value = "CD56-PC7"
$ws.Range("N7").Value = "CD33-PC7"
$ws.Range("N8").Value = "CD38-PC7"
$ws.Range("N9").Value = "CD10-PC7"

# Remaining cells: Tube numbers (G3:G10) and repeated FSC/SSC/CD45-ECD labels
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = "FSC"
$ws.Range("I3").Value = "SSC"

$ws.Range("G4").Value = 2
$ws.Range("H4").Value = "FSC"
$ws.Range("I4").Value = "SSC"
$ws.Range("L4").Value = "CD45-ECD"

$ws.Range("G5").Value = 3
$ws.Range("H5").Value = "FSC"
$ws.Range("I5").Value = "SSC"
$ws.Range("L5").Value = "CD45-ECD"

$ws.Range("G6").Value = 4
$ws.Range("H6").Value = "FSC"
$ws.Range("I6").Value = "SSC"
$ws.Range("L6").Value = "CD45-ECD"

$ws.Range("G7").Value = 5
$ws.Range("H7").Value = "FSC"
$ws.Range("I7").Value = "SSC"
$ws.Range("L7").Value = "CD45-ECD"

$ws.Range("G8").Value = 6
$ws.Range("H8").Value = "FSC"
$ws.Range("I8").Value = "SSC"
$ws.Range("L8").Value = "CD45-ECD"

$ws.Range("G9").Value = 7
$ws.Range("H9").Value = "FSC"
$ws.Range("I9").Value = "SSC"
$ws.Range("L9").Value = "CD45-ECD"

$ws.Range("G10").Value = 8
$ws.Range("H10").Value = "FSC"
$ws.Range("I10").Value = "SSC"

# --- Formatting ---
# Bold header row for the new table
$ws.Range("G2:N2").Font.Bold = $true

# Column widths (engine stores ColumnWidth + 0.8333... as the saved <col> width,
# so we back the desired values out through that offset)
$ws.Columns.Item(9).ColumnWidth = 3.3333333333333335
$ws.Columns.Item(10).ColumnWidth = 11.166666666666666
$ws.Columns.Item(11).ColumnWidth = 9.833333333333334
$ws.Columns.Item(12).ColumnWidth = 8.666666666666666
$ws.Columns.Item(13).ColumnWidth = 8.666666666666666

# --- View state ---
$ws.Activate()
$ws.Range("L10").Select()
$excel.ActiveWindow.Zoom = 145

# --- Page setup ---
$ws.PageSetup.Orientation = 1
